# Fixed bug by renaming transect_chl_final to api_chl_final
# (In this workbook, the erroneous label was the shared string "attribute"
# used as the header of the first column on the CategoricalVariables sheet;
# it should read "attributeName" to match the header used on the
# ColumnHeaders / StationData sheets.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CategoricalVariables")
$ws.Range("A1").Value = "attributeName"
